$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.242.40"
$ws.Range("E2").Value = "  +1.98%  "

$ws.Range("D3").Value = "3.616.96"
$ws.Range("E3").Value = "  +6.02%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'237.70"
$ws.Range("E5").Value = "  +1.79%  "

$ws.Range("D6").Value = "'659.00"
$ws.Range("E6").Value = "  +6.27%  "

$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("E8").Value = "  +3.49%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "'0.993"
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("D11").Value = "3.615.62"
$ws.Range("E11").Value = "  +6.02%  "

$ws.Range("D12").Value = "'42.54"
$ws.Range("E12").Value = "  -1.83%  "

$ws.Range("D13").Value = "'0.201"
$ws.Range("E13").Value = "  +0.85%  "

$ws.Range("D14").Value = "'6.30"
$ws.Range("E14").Value = "  +0.31%  "

$ws.Range("D15").Value = "4.294.38"
$ws.Range("E15").Value = "  +5.99%  "

$ws.Range("D16").Value = "95.427.09"
$ws.Range("E16").Value = "  +2.43%  "

$ws.Range("E17").Value = "  +3.05%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.622.38"
$ws.Range("E18").Value = "  +6.23%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'13.09"
$ws.Range("E19").Value = "  +12.54%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'7.88"
$ws.Range("E20").Value = "  -4.74%  "

$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("E22").Value = "  +5.76%  "

$ws.Range("E23").Value = "  -5.60%  "

$ws.Range("D24").Value = "'505.91"
$ws.Range("E24").Value = "  +1.54%  "

$ws.Range("E25").Value = "  +7.52%  "

$ws.Range("D26").Value = "'6.61"
$ws.Range("E26").Value = "  -2.59%  "

$ws.Range("D27").Value = "'95.81"
$ws.Range("E27").Value = "  +6.38%  "

$ws.Range("D28").Value = "'12.62"
$ws.Range("E28").Value = "  +4.93%  "

$ws.Range("D29").Value = "3.816.13"
$ws.Range("E29").Value = "  +6.32%  "

$ws.Range("D30").Value = "'3.18"
$ws.Range("E30").Value = "  +17.01%  "

$ws.Range("D31").Value = "'11.33"
$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("E35").Value = "  +1.27%  "

$ws.Range("D36").Value = "'32.00"
$ws.Range("E36").Value = "  +10.60%  "

$ws.Range("E37").Value = "  +1.45%  "

$ws.Range("D38").Value = "'575.09"
$ws.Range("E38").Value = "  +3.20%  "

$ws.Range("D39").Value = "'8.14"
$ws.Range("E39").Value = "  +8.76%  "

$ws.Range("E40").Value = "  +5.22%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").Value = "'0.923"
$ws.Range("E42").Value = "  +3.26%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Value = "'34.92"
$ws.Range("E44").Value = "  +46.04%  "

$ws.Range("E45").Value = "  +1.59%  "

$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").Value = "'5.66"
$ws.Range("E47").Value = "  +3.68%  "

$ws.Range("E48").Value = "  +6.02%  "

$ws.Range("E49").Value = "  -1.03%  "

$ws.Range("E50").Value = "  -3.36%  "

$ws.Range("D51").Value = "'53.55"
$ws.Range("E51").Value = "  +0.95%  "

